$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SQL queries (join conditions now use the explicit *_id columns
# instead of the generic "id" column, and a trailing space after WHERE was
# trimmed in all but the Treatment query).
$ws.Range("C2").Value = "SELECT `r`n    COUNT(DISTINCT dgn.diagnosis) AS Diagnoses,`r`n    COUNT(DISTINCT prt.participant_id) AS Participants,`r`n    COUNT(DISTINCT std.study_id) AS Studies`r`nFROM `r`n    df_study std`r`nLEFT JOIN `r`n    df_participant prt ON std.study_id = prt.`"study.study_id`"`r`nLEFT JOIN `r`n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`r`nWHERE`r`n    std.dbgap_accession = 'phs002529' AND dgn.diagnosis= '8800/3 : Sarcoma, NOS';"
$ws.Range("B2").Value = "SELECT DISTINCT`r`n     std.dbgap_accession AS `"dbGaP Accession`",`r`n     std.study_name AS `"Study Name`"`r`nFROM `r`n    df_study std`r`nLEFT JOIN `r`n    df_participant prt ON std.study_id = prt.`"study.study_id`"`r`nLEFT JOIN `r`n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`r`nWHERE`r`n    std.dbgap_accession = 'phs002529' AND dgn.diagnosis= '8800/3 : Sarcoma, NOS'"
$ws.Range("B3").Value = "SELECT DISTINCT`r`n    prt.participant_id AS `"Participant Id`",`r`n    prt.race AS `"Race`",`r`n    prt.sex_at_birth AS `"Sex at Birth`",`r`n    std.dbgap_accession AS `"dbGaP Accession`"`r`nFROM `r`n    df_study std`r`nLEFT JOIN `r`n    df_participant prt ON std.study_id = prt.`"study.study_id`"`r`nLEFT JOIN `r`n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`r`nWHERE`r`n    std.dbgap_accession = 'phs002529' AND dgn.diagnosis= '8800/3 : Sarcoma, NOS'`r`nORDER BY `r`n    prt.participant_id ASC`r`nLIMIT 100;"
$ws.Range("B4").Value = "SELECT DISTINCT`r`n    prt.participant_id AS `"Participant Id`",`r`n    dgn.diagnosis_id AS `"Diagnosis Id`",`r`n    dgn.diagnosis AS `"Diagnosis`",`r`n    dgn.diagnosis_classification_system AS `"Diagnosis Classification System`",`r`n    dgn.diagnosis_basis AS `"Diagnosis Basis`",`r`n    dgn.tumor_classification AS `"Tumor Classification`",`r`n    dgn.anatomic_site AS `"Anatomic Site`",`r`n    CASE `r`n    WHEN dgn.age_at_diagnosis = -999 THEN 'Not Reported'`r`n    WHEN dgn.age_at_diagnosis >= 1000 THEN `r`n        substr(dgn.age_at_diagnosis, 1, length(dgn.age_at_diagnosis) - 3) || ',' || substr(dgn.age_at_diagnosis, -3)`r`n    ELSE `r`n        dgn.age_at_diagnosis `r`nEND AS `"Age at Diagnosis (days)`",`r`n    std.dbgap_accession AS `"dbGaP Accession`"`r`nFROM `r`n    df_study std`r`nLEFT JOIN `r`n    df_participant prt ON std.study_id = prt.`"study.study_id`"`r`nLEFT JOIN `r`n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`r`nWHERE`r`n    std.dbgap_accession = 'phs002529' AND dgn.diagnosis= '8800/3 : Sarcoma, NOS'`r`nORDER BY `r`n    prt.participant_id ASC`r`nLIMIT 100;"
$ws.Range("B5").Value = "SELECT`r`n    DISTINCT prt.participant_id AS `"Participant Id`",`r`n    trt.treatment_id AS `"Treatment Id`",`r`n    CASE `r`n    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'`r`n    WHEN trt.age_at_treatment_start >= 1000 THEN `r`n        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)`r`n    ELSE `r`n        trt.age_at_treatment_start `r`nEND AS `"Age at Treatment Start`",`r`n    CASE `r`n    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'`r`n    WHEN trt.age_at_treatment_end >= 1000 THEN `r`n        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)`r`n    ELSE `r`n        trt.age_at_treatment_end `r`nEND AS `"Age at Treatment End`",`r`n    trt.treatment_type AS `"Treatment Type`",`r`n    CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS `"Treatment Agent`",`r`n    std.dbgap_accession AS `"dbGaP Accession`"`r`nFROM `r`n    df_study std`r`nLEFT JOIN `r`n    df_participant prt ON std.study_id = prt.`"study.study_id`"`r`nLEFT JOIN `r`n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`r`nWHERE `r`n    std.dbgap_accession = 'phs002529' AND dgn.diagnosis= '8800/3 : Sarcoma, NOS'`r`nORDER BY `r`n    trt.treatment_id ASC`r`nLIMIT 100;"
$ws.Range("B6").Value = "SELECT DISTINCT`r`n    prt.participant_id AS `"Participant Id`",`r`n    trr.treatment_response_id AS `"Treatment Response Id`",`r`n    trr.response AS `"Response`",`r`n    CASE `r`n        WHEN trr.age_at_response = -999 THEN 'Not Reported'`r`n        WHEN trr.age_at_response >= 1000 THEN `r`n            substr(trr.age_at_response, 1, length(trr.age_at_response) - 3) || ',' || substr(trr.age_at_response, -3)`r`n        ELSE `r`n            trr.age_at_response `r`n    END AS `"Age at Response`",`r`n    trr.response_category AS `"Response Category`",`r`n    trr.response_system AS `"Response System`",`r`n    std.dbgap_accession AS `"dbGaP Accession`"`r`nFROM `r`n    df_study std`r`nLEFT JOIN `r`n    df_participant prt ON std.study_id = prt.`"study.study_id`"`r`nLEFT JOIN `r`n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`r`nWHERE`r`n    std.dbgap_accession = 'phs002529' AND dgn.diagnosis= '8800/3 : Sarcoma, NOS'`r`nORDER BY `r`n    prt.participant_id ASC`r`nLIMIT 100;"
$ws.Range("B7").Value = "SELECT DISTINCT`r`n    prt.participant_id AS `"Participant Id`",`r`n    srv.survival_id AS `"Survival Id`",`r`n    srv.last_known_survival_status AS `"Last Known Survival Status`",`r`n    CASE `r`n    WHEN srv.age_at_last_known_survival_status = -999 THEN 'Not Reported'`r`n    WHEN srv.age_at_last_known_survival_status >= 1000 THEN `r`n        substr(srv.age_at_last_known_survival_status, 1, length(srv.age_at_last_known_survival_status) - 3) || ',' || substr(srv.age_at_last_known_survival_status, -3)`r`n    ELSE `r`n        srv.age_at_last_known_survival_status `r`nEND AS `"Age at Last Known Survival Status`",`r`n    srv.first_event AS `"First Event`",`r`n    srv.cause_of_death AS `"Cause of Death`",`r`n    std.dbgap_accession AS `"dbGaP Accession`"`r`nFROM `r`n    df_study std`r`nLEFT JOIN `r`n    df_participant prt ON std.study_id = prt.`"study.study_id`"`r`nLEFT JOIN `r`n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`r`nLEFT JOIN `r`n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`r`nWHERE`r`n    std.dbgap_accession = 'phs002529' AND dgn.diagnosis= '8800/3 : Sarcoma, NOS'`r`nORDER BY `r`n    prt.participant_id ASC`r`nLIMIT 100;"

# Update sheet view: scroll back to top-left (drop the saved topLeftCell)
# and move the active selection to B2.
$null = $ws.Range("B2").Select()
